$d = $word.ActiveDocument

# The "_GoBack" bookmark marks the position of the author's last edit. In
# the target revision it ends up (still zero-width) right after the new,
# final paragraph "For få møder i starten af forløbet". Word COM only
# reliably "carries" a bookmark forward when plain text is inserted via
# the bookmark's own (live, non-duplicated) Range - paragraph breaks
# inserted that way get left behind. So: insert all of the new text as one
# flat run (with unique placeholder markers where paragraph breaks belong)
# straight after the bookmark, then turn each placeholder into a real
# paragraph break via Find/Replace - which still keeps the bookmark
# anchored right after the inserted text.

$newText = (
    "@@SPLIT1@@" +
    "For abstrakte aftaler ved start, skulle have været mere konkrete for at undgå misforståelser" +
    "@@SPLIT2@@" +
    "For få møder i starten af forløbet"
)

if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Range.InsertAfter($newText)
} else {
    # Fallback: no _GoBack bookmark present - just tack the new paragraphs
    # on after the last paragraph of body text.
    $lastPara = $d.Paragraphs.Item(8)
    $lastPara.Range.InsertAfter($newText)
}

$d.Content.Find.Execute("@@SPLIT1@@", $true, $false, $false, $false, $false, $true, 1, $false, "^p", 2) | Out-Null
$d.Content.Find.Execute("@@SPLIT2@@", $true, $false, $false, $false, $false, $true, 1, $false, "^p", 2) | Out-Null
